{"js": "const replacements = [\n  [\"2024-03-20 Wednesday\", \"2024-03-21 Thursday\"],\n  [\"438\u00d79=3942\", \"463\u00d73=1389\"],\n  [\"958\u00d73=2874\", \"737\u00d77=5159\"],\n  [\"193\u00d73=579\", \"688\u00d79=6192\"],\n  [\"108\u00d75=540\", \"395\u00d76=2370\"],\n  [\"840\u00d77=5880\", \"960\u00d78=7680\"],\n  [\"150\u00d79=1350\", \"331\u00d78=2648\"],\n  [\"519\u00d78=4152\", \"173\u00d79=1557\"],\n  [\"952\u00d73=2856\", \"413\u00d76=2478\"],\n  [\"921\u00d76=5526\", \"403\u00d74=1612\"],\n  [\"204\u00d77=1428\", \"126\u00d79=1134\"],\n  [\"448\u00d78=3584\", \"351\u00d75=1755\"],\n  [\"625\u00d72=1250\", \"642\u00d73=1926\"],\n  [\"128\u00d73=384\", \"293\u00d75=1465\"],\n  [\"119\u00d75=595\", \"126\u00d74=504\"],\n  [\"267\u00d74=1068\", \"298\u00d78=2384\"],\n  [\"199\u00d78=1592\", \"320\u00d78=2560\"],\n  [\"659\u00d72=1318\", \"784\u00d72=1568\"],\n  [\"796\u00d72=1592\", \"877\u00d76=5262\"],\n  [\"203\u00d76=1218\", \"767\u00d75=3835\"],\n  [\"164\u00d78=1312\", \"300\u00d75=1500\"],\n  [\"520\u00d73=1560\", \"733\u00d74=2932\"],\n  [\"641\u00d75=3205\", \"421\u00d77=2947\"],\n  [\"479\u00d73=1437\", \"683\u00d79=6147\"],\n  [\"648\u00d73=1944\", \"135\u00d72=270\"],\n  [\"339\u00d76=2034\", \"478\u00d76=2868\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + before);\n  }\n\n  for (const range of results.items) {\n    range.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Before = \"2024-03-20 Wednesday\"; After = \"2024-03-21 Thursday\" },\n    @{ Before = \"438\u00d79=3942\"; After = \"463\u00d73=1389\" },\n    @{ Before = \"958\u00d73=2874\"; After = \"737\u00d77=5159\" },\n    @{ Before = \"193\u00d73=579\"; After = \"688\u00d79=6192\" },\n    @{ Before = \"108\u00d75=540\"; After = \"395\u00d76=2370\" },\n    @{ Before = \"840\u00d77=5880\"; After = \"960\u00d78=7680\" },\n    @{ Before = \"150\u00d79=1350\"; After = \"331\u00d78=2648\" },\n    @{ Before = \"519\u00d78=4152\"; After = \"173\u00d79=1557\" },\n    @{ Before = \"952\u00d73=2856\"; After = \"413\u00d76=2478\" },\n    @{ Before = \"921\u00d76=5526\"; After = \"403\u00d74=1612\" },\n    @{ Before = \"204\u00d77=1428\"; After = \"126\u00d79=1134\" },\n    @{ Before = \"448\u00d78=3584\"; After = \"351\u00d75=1755\" },\n    @{ Before = \"625\u00d72=1250\"; After = \"642\u00d73=1926\" },\n    @{ Before = \"128\u00d73=384\"; After = \"293\u00d75=1465\" },\n    @{ Before = \"119\u00d75=595\"; After = \"126\u00d74=504\" },\n    @{ Before = \"267\u00d74=1068\"; After = \"298\u00d78=2384\" },\n    @{ Before = \"199\u00d78=1592\"; After = \"320\u00d78=2560\" },\n    @{ Before = \"659\u00d72=1318\"; After = \"784\u00d72=1568\" },\n    @{ Before = \"796\u00d72=1592\"; After = \"877\u00d76=5262\" },\n    @{ Before = \"203\u00d76=1218\"; After = \"767\u00d75=3835\" },\n    @{ Before = \"164\u00d78=1312\"; After = \"300\u00d75=1500\" },\n    @{ Before = \"520\u00d73=1560\"; After = \"733\u00d74=2932\" },\n    @{ Before = \"641\u00d75=3205\"; After = \"421\u00d77=2947\" },\n    @{ Before = \"479\u00d73=1437\"; After = \"683\u00d79=6147\" },\n    @{ Before = \"648\u00d73=1944\"; After = \"135\u00d72=270\" },\n    @{ Before = \"339\u00d76=2034\"; After = \"478\u00d76=2868\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Before\n    $find.Replacement.Text = $r.After\n    $find.Execute([ref]$r.Before, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, 1, [ref]$false, [ref]$r.After, 2)\n}\n"}
